# Applies the "I" -> "We"/"we" voice change to several paragraphs in the
# "Main results from test runs" document, and removes the stray _GoBack
# bookmark at the very end of the document, as described by the commit.
#
# NOTE: We intentionally avoid passing the replacement text through
# Find.Execute's own "ReplaceWith" argument when it contains an apostrophe,
# because this runtime auto-"smart-quotes" straight apostrophes inserted
# that way. Instead we locate the text with Find.Execute (no replacement),
# then set .Text directly on the matched Range, which preserves the
# original straight-apostrophe style used throughout the document.

$d = $word.ActiveDocument

function Replace-Text($doc, $find, $replace) {
    $r = $doc.Content
    $found = $r.Find.Execute($find, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: text not found:" $find
        return
    }
    $r.Text = $replace
}

# Paragraph: "I found that all the time after the crash, ... ntdll!RtlUserThreadStart, which as I noticed ..."
Replace-Text $d "I found that all the time after the crash" "We found that all the time after the crash"
Replace-Text $d "In addition, I found that the specific system call" "In addition, We found that the specific system call"
Replace-Text $d ", which as I noticed after reading the forums" ", which as we noticed after reading the forums"

# Paragraph: "I couldn't get a superficial look ... Looking at memory usage as well as disk operations, I can assume ..."
Replace-Text $d "I couldn't get a superficial look at what caused the problem" "We couldn't get a superficial look at what caused the problem"
Replace-Text $d "Looking at memory usage as well as disk operations, I can assume" "Looking at memory usage as well as disk operations, We can assume"

# Paragraph: "Also complete and analyse the same testing for application with DB data source. DB that I was using is SQLite."
Replace-Text $d "DB that I was using is SQLite" "DB that We was using is SQLite"

# Paragraph: "I'm sure that the performance when working with a database depends heavily on the database itself."
Replace-Text $d "I'm sure that the performance when working with a database" "We are sure that the performance when working with a database"

# Paragraph: "All the queries that had errors here, too, are comments and dashboard queries. Since they are write and read from disk, I think that the problem is with disk operations and queues."
Replace-Text $d "Since they are write and read from disk, I think that the problem" "Since they are write and read from disk, We think that the problem"

# Paragraph: "In conclusion, I would like to say that the application by modern standards is quite productive ..."
Replace-Text $d "In conclusion, I would like to say that the application" "In conclusion, we would like to say that the application"

# Remove the trailing _GoBack bookmark (present at the very end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}
